# FN-1945: change fees paid report header
# "Monthly fees paid to UKEF" -> "Fees paid to UKEF for the period"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("invalid-utilisation-report")

$ws.Range("H1").Value = "Fees paid to UKEF for the period"
